# -----------------------------------------------------------------------
# Reproduces the authors edit:
#   * adds a new shared string "hello@tide.com"
#   * changes the value of C5 (on "Tide test data") from
#     "hellothere@tide.com" to "hello@tide.com", while keeping its
#     existing (quote-prefixed hyperlink) cell style untouched
#   * keeps C5's mailto hyperlink pointing at hellothere@tide.com but
#     gives it an explicit display text of "hellothere@tide.com" (so the
#     hyperlink text differs from the new cell text)
#   * moves the worksheet selection from A10:D13 to C6
#   * best-effort: repositions/resizes the workbook window
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet's Hyperlinks collection in this runtime can only be
# cleared/rebuilt as a whole -- individual Hyperlink.Delete() calls are a
# no-op, so remove every hyperlink and recreate them all in their
# original order/addresses. Only the C5 entry gets a custom
# TextToDisplay (which is what produces the new `display="..."`
# attribute in the saved XML).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:hellothere@tide.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:hellothere@tide.com", [Type]::Missing, [Type]::Missing, "hellothere@tide.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:P@ss1234")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:P@ss1234")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:jackjone@tide.com")

# Adding a hyperlink re-writes the anchor cell's text/style (it applies a
# fresh "hyperlink" font style). Re-assign each cell's value (leading
# apostrophe = keep it stored/text-only, matching the workbook's
# original quote-prefixed hyperlink style) so every cell ends up with
# its original text and original style, except C5 which gets the new
# address text.
$ws.Range("C3").Value = "'hellothere@tide.com"
$ws.Range("C4").Value = "'P@ss1234"
$ws.Range("C5").Value = "'hello@tide.com"
$ws.Range("C6").Value = "'P@ss1234"
$ws.Range("C12").Value = "'jackjone@tide.com"

# Move the active selection to C6.
$ws.Range("C6").Select()

# Best-effort: match the author's Excel window position/size.
$win = $excel.ActiveWindow
$win.Left = 1400
$win.Top = 1400
$win.Width = 14400
$win.Height = 7360

Write-Host "edit.ps1 complete"
